$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.127.13"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "1.651.13"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5204"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2661"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06320"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.431"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.649.50"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "1.880.70"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5458"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("D16").Value = "0.0₅8217"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "26.158.60"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.666"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.090"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1239"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.217"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.426"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06024"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.281"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.319"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9789"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.771"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5924"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01589"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.956"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8649"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "1.041.95"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "1.794.67"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.464"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.82%  "
